$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Week 1 (rows 7-11): hours halved/adjusted and activity notes trimmed
# (the recurring ", desarrollo proyecto AT04" suffix was dropped from the log)
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "Construccion y envio de 500MD, AT37, AT38"

$ws.Range("E8").Value = 1.5
$ws.Range("F8").Value = "Construccion y envio de 500MD "

$ws.Range("E9").Value = 1.5
$ws.Range("F9").Value = "Construccion y envio de 500MD, completar plantilla EUC."

$ws.Range("E11").Value = 2.5
$ws.Range("F11").Value = "Construccion y envio de 500MD, meeting con equipo de Mexico sobre proyecto 500MD, creacion de mapping para el proyecto 500MD"

# Week 2 (rows 12-16)
$ws.Range("E12").Value = 2.5
$ws.Range("F12").Value = "Construccion y envio de 500MD, AT37"

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = "Construccion y envio de 500MD, AT38, tablas manuales auxiliar contable, y carga de las mismas al portal del LRR"

$ws.Range("E14").Value = 0.5

$ws.Range("E15").Value = 0.5

# Move the (stale) selection that was left on the sheet
$ws.Range("D20").Select()
